# Apply BoM updates: two component groups lost references (C11/C28 removed from
# the "100n" capacitor group, D2/D6 removed from the SK6812 RGB LED group),
# which reduces the total/fitted component counts by 4 in the summary header.

$wb = $excel.ActiveWorkbook

# --- BoM sheet ---
$bom = $wb.Worksheets.Item("BoM")

# Row 4 ("100n" small capacitors): drop C11 and C28 from References, update qty.
$bom.Range("D12").Value = "C1 C12 C13 C14 C16 C18 C19 C20 C21 C22 C25 C27 C36 C39"
$bom.Range("G12").Value = "14"

# Row 11 (SK6812 RGB LED): drop D2 and D6 from References, update qty.
$bom.Range("D19").Value = "D3 D4"
$bom.Range("G19").Value = "2"

# Header summary counts.
$bom.Range("F3").Value = "112 (85 SMD/ 25 THT)"
$bom.Range("F4").Value = "102 (82 SMD/ 20 THT)"
$bom.Range("F6").Value = 102

# Slightly narrower References column.
$bom.Columns.Item(4).ColumnWidth = 59.7109375

# --- DNF sheet mirrors the same summary header ---
$dnf = $wb.Worksheets.Item("DNF")
$dnf.Range("F3").Value = "112 (85 SMD/ 25 THT)"
$dnf.Range("F4").Value = "102 (82 SMD/ 20 THT)"
$dnf.Range("F6").Value = 102
